$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "take foreign spicules"
$ws.Range("A2").Value = $null
$ws.Range("A4").Value = $null
